$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.7321483333333333
$ws.Range("H2").Value = 2.196445
$ws.Range("I2").Value = 0.05113520435363902
$ws.Range("J2").Value = 0.05113520435363902
$ws.Range("M2").Value = 1.442875
$ws.Range("N2").Value = 4.328625
$ws.Range("O2").Value = 0.02047893724893121
$ws.Range("P2").Value = 0.02047893724893121
$ws.Range("Q2").Value = 1.056398526458333
$ws.Range("R2").Value = 9.507586738124999
$ws.Range("S2").Value = 0.001047194641169447
$ws.Range("T2").Value = 0.001047194641169447
$ws.Range("G3").Value = 0.7321483333333333
$ws.Range("H3").Value = 2.196445
$ws.Range("I3").Value = 0.05113520435363902
$ws.Range("J3").Value = 0.05113520435363902
$ws.Range("O3").Value = 0.1473796107804731
$ws.Range("P3").Value = 0.1473796107804731
$ws.Range("Q3").Value = 7.602523596121666
$ws.Range("R3").Value = 68.42271236509499
$ws.Range("S3").Value = 0.007536286514819272
$ws.Range("T3").Value = 0.00753628651481927
$ws.Range("G4").Value = 0.7321483333333333
$ws.Range("H4").Value = 2.196445
$ws.Range("I4").Value = 0.05113520435363902
$ws.Range("J4").Value = 0.05113520435363902
$ws.Range("M4").Value = 27.934719
$ws.Range("N4").Value = 83.804157
$ws.Range("O4").Value = 0.3964815784233052
$ws.Range("P4").Value = 0.3964815784233051
$ws.Range("Q4").Value = 20.452357957985
$ws.Range("R4").Value = 184.071221621865
$ws.Range("S4").Value = 0.02027416653512907
$ws.Range("T4").Value = 0.02027416653512906
$ws.Range("G5").Value = 0.7321483333333333
$ws.Range("H5").Value = 2.196445
$ws.Range("I5").Value = 0.05113520435363902
$ws.Range("J5").Value = 0.05113520435363902
$ws.Range("M5").Value = 30.695086
$ws.Range("N5").Value = 92.085258
$ws.Range("O5").Value = 0.4356598735472906
$ws.Range("P5").Value = 0.4356598735472905
$ws.Range("Q5").Value = 22.47335605642333
$ws.Range("R5").Value = 202.26020450781
$ws.Range("S5").Value = 0.02227755666252124
$ws.Range("T5").Value = 0.02227755666252123
$ws.Range("I6").Value = 0.7165747117895102
$ws.Range("J6").Value = 0.7165747117895102
$ws.Range("M6").Value = 1.442875
$ws.Range("N6").Value = 4.328625
$ws.Range("O6").Value = 0.02047893724893121
$ws.Range("P6").Value = 0.02047893724893121
$ws.Range("Q6").Value = 14.80366567808333
$ws.Range("R6").Value = 133.23299110275
$ws.Range("S6").Value = 0.01467468855690834
$ws.Range("T6").Value = 0.01467468855690834
$ws.Range("I7").Value = 0.7165747117895102
$ws.Range("J7").Value = 0.7165747117895102
$ws.Range("O7").Value = 0.1473796107804731
$ws.Range("P7").Value = 0.1473796107804731
$ws.Range("R7").Value = 958.8303403227781
$ws.Range("S7").Value = 0.1056085021186677
$ws.Range("T7").Value = 0.1056085021186677
$ws.Range("I8").Value = 0.7165747117895102
$ws.Range("J8").Value = 0.7165747117895102
$ws.Range("M8").Value = 27.934719
$ws.Range("N8").Value = 83.804157
$ws.Range("O8").Value = 0.3964815784233052
$ws.Range("P8").Value = 0.3964815784233051
$ws.Range("Q8").Value = 286.605728761814
$ws.Range("R8").Value = 2579.451558856326
$ws.Range("S8").Value = 0.28410867278853
$ws.Range("T8").Value = 0.2841086727885299
$ws.Range("I9").Value = 0.7165747117895102
$ws.Range("J9").Value = 0.7165747117895102
$ws.Range("M9").Value = 30.695086
$ws.Range("N9").Value = 92.085258
$ws.Range("O9").Value = 0.4356598735472906
$ws.Range("P9").Value = 0.4356598735472905
$ws.Range("Q9").Value = 314.9266506828494
$ws.Range("R9").Value = 2834.339856145644
$ws.Range("S9").Value = 0.3121828483254042
$ws.Range("T9").Value = 0.3121828483254042
$ws.Range("G10").Value = 2.568000333333333
$ws.Range("H10").Value = 7.704001
$ws.Range("I10").Value = 0.1793560346266988
$ws.Range("J10").Value = 0.1793560346266988
$ws.Range("M10").Value = 1.442875
$ws.Range("N10").Value = 4.328625
$ws.Range("O10").Value = 0.02047893724893121
$ws.Range("P10").Value = 0.02047893724893121
$ws.Range("Q10").Value = 3.705303480958333
$ws.Range("R10").Value = 33.347731328625
$ws.Range("S10").Value = 0.003673020978337297
$ws.Range("T10").Value = 0.003673020978337297
$ws.Range("G11").Value = 2.568000333333333
$ws.Range("H11").Value = 7.704001
$ws.Range("I11").Value = 0.1793560346266988
$ws.Range("J11").Value = 0.1793560346266988
$ws.Range("O11").Value = 0.1473796107804731
$ws.Range("P11").Value = 0.1473796107804731
$ws.Range("Q11").Value = 26.66574823728567
$ws.Range("R11").Value = 239.991734135571
$ws.Range("S11").Value = 0.02643342257441192
$ws.Range("T11").Value = 0.02643342257441192
$ws.Range("G12").Value = 2.568000333333333
$ws.Range("H12").Value = 7.704001
$ws.Range("I12").Value = 0.1793560346266988
$ws.Range("J12").Value = 0.1793560346266988
$ws.Range("M12").Value = 27.934719
$ws.Range("N12").Value = 83.804157
$ws.Range("O12").Value = 0.3964815784233052
$ws.Range("P12").Value = 0.3964815784233051
$ws.Range("Q12").Value = 71.73636770357301
$ws.Range("R12").Value = 645.627309332157
$ws.Range("S12").Value = 0.0711113637085385
$ws.Range("T12").Value = 0.0711113637085385
$ws.Range("G13").Value = 2.568000333333333
$ws.Range("H13").Value = 7.704001
$ws.Range("I13").Value = 0.1793560346266988
$ws.Range("J13").Value = 0.1793560346266988
$ws.Range("M13").Value = 30.695086
$ws.Range("N13").Value = 92.085258
$ws.Range("O13").Value = 0.4356598735472906
$ws.Range("P13").Value = 0.4356598735472905
$ws.Range("Q13").Value = 78.82499107969534
$ws.Range("R13").Value = 709.4249197172579
$ws.Range("S13").Value = 0.07813822736541105
$ws.Range("T13").Value = 0.07813822736541105
$ws.Range("G14").Value = 0.7579039999999999
$ws.Range("H14").Value = 2.273712
$ws.Range("I14").Value = 0.05293404923015203
$ws.Range("J14").Value = 0.05293404923015203
$ws.Range("M14").Value = 1.442875
$ws.Range("N14").Value = 4.328625
$ws.Range("O14").Value = 0.02047893724893121
$ws.Range("P14").Value = 0.02047893724893121
$ws.Range("Q14").Value = 1.093560734
$ws.Range("R14").Value = 9.842046605999998
$ws.Range("S14").Value = 0.001084033072516119
$ws.Range("T14").Value = 0.001084033072516119
$ws.Range("G15").Value = 0.7579039999999999
$ws.Range("H15").Value = 2.273712
$ws.Range("I15").Value = 0.05293404923015203
$ws.Range("J15").Value = 0.05293404923015203
$ws.Range("O15").Value = 0.1473796107804731
$ws.Range("P15").Value = 0.1473796107804731
$ws.Range("Q15").Value = 7.869966755728
$ws.Range("R15").Value = 70.82970080155199
$ws.Range("S15").Value = 0.007801399572574207
$ws.Range("T15").Value = 0.007801399572574206
$ws.Range("G16").Value = 0.7579039999999999
$ws.Range("H16").Value = 2.273712
$ws.Range("I16").Value = 0.05293404923015203
$ws.Range("J16").Value = 0.05293404923015203
$ws.Range("M16").Value = 27.934719
$ws.Range("N16").Value = 83.804157
$ws.Range("O16").Value = 0.3964815784233052
$ws.Range("P16").Value = 0.3964815784233051
$ws.Range("Q16").Value = 21.171835268976
$ws.Range("R16").Value = 190.546517420784
$ws.Range("S16").Value = 0.02098737539110762
$ws.Range("T16").Value = 0.02098737539110761
$ws.Range("G17").Value = 0.7579039999999999
$ws.Range("H17").Value = 2.273712
$ws.Range("I17").Value = 0.05293404923015203
$ws.Range("J17").Value = 0.05293404923015203
$ws.Range("M17").Value = 30.695086
$ws.Range("N17").Value = 92.085258
$ws.Range("O17").Value = 0.4356598735472906
$ws.Range("P17").Value = 0.4356598735472905
$ws.Range("Q17").Value = 23.263928459744
$ws.Range("R17").Value = 209.375356137696
$ws.Range("S17").Value = 0.02306124119395409
$ws.Range("T17").Value = 0.02306124119395408
